$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BDMFL")
$ws.Range("B2").Value = 1
